$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "municipio-nombre" column (D) is re-curated from a measure to a
# dimension, matching the treatment already used by provincia-nombre (E)
# and comarca-nombre (K):
#   D2: iaest-measure:municipio-nombre -> sdmx-dimension:refArea
#   D3: medida                         -> dim
#   D4: xsd:int                        -> URI-Municipio
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"
